# Generate Report for Handoff
#
# The workbook tracks localization handoff status for a set of files.
# The fixture data changes from 3 files (2 images + 1 markdown, with a
# dependency relationship) to 4 markdown files: two "callee" files that
# are referenced (Included) by two "caller" files that depend on them.
#
# This script rewrites the "Overview" sheet and the per-locale ("zh-cn",
# "de-de") detail sheets to reflect the new fixture set, and appends a
# 4th data row (row 5) to every sheet for the new "callerMd2.md" file.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "Overview": columns A..D -> File Name, zh-cn, de-de, Latest Handoff Date
# ---------------------------------------------------------------------
$ov = $wb.Worksheets.Item("Overview")

$ovNames = @("calleeMd1.md", "calleeMd2.md", "callerMd1.md", "callerMd2.md")
$ovLinks = @(
    "https://github.com/OpenLocalizationTest/oltest/blob/6cb7b01ae03b28ec25e20f4b9bba2f4dadeffb1d/e2e/calleeMd1.md",
    "https://github.com/OpenLocalizationTest/oltest/blob/6cb7b01ae03b28ec25e20f4b9bba2f4dadeffb1d/e2e/calleeMd2.md",
    "https://github.com/OpenLocalizationTest/oltest/blob/6cb7b01ae03b28ec25e20f4b9bba2f4dadeffb1d/e2e/callerMd1.md",
    "https://github.com/OpenLocalizationTest/oltest/blob/6cb7b01ae03b28ec25e20f4b9bba2f4dadeffb1d/e2e/callerMd2.md"
)

# Hyperlinks.Delete() wipes every hyperlink on the sheet (not just the
# ones in the calling Range), so clear once up front and re-Add in the
# final row order rather than trying to patch the 3 existing links.
$ov.Hyperlinks.Delete()

for ($i = 0; $i -lt 4; $i++) {
    $row = $i + 2
    $ov.Range("B$row").Value = "Ready for handoff"
    $ov.Range("C$row").Value = "Ready for handoff"
    $ov.Range("D$row").Value = "2016-03-19 17:15:08"
    $ov.Hyperlinks.Add($ov.Range("A$row"), $ovLinks[$i], $null, $null, $ovNames[$i]) | Out-Null
}

# ---------------------------------------------------------------------
# Shared layout for the "zh-cn" / "de-de" detail sheets (columns A..L):
#   A Source File Name   B File Extension   C Status
#   D Latest Handoff File   E Latest Handoff Datetime   F Latest Target File
#   G Latest Handback File  H Latest Handback DateTime  I Reference Tokens
#   J Handoff Reason        K Dependency From            L Error Detail
# ---------------------------------------------------------------------
function Update-DetailSheet($ws, $xlfSuffix, $handoffLink, $xlfLinkBase, $handoffDatetime) {
    $names = @("calleeMd1.md", "calleeMd2.md", "callerMd1.md", "callerMd2.md")
    $xlfHashes = @(
        "e8f5ecec2b522eb147a4ff0ca19ca72e17f2186d",
        "63b76063f058ecc63ff1dda71ea2a67db72ae6e1",
        "a3bf9f4e7fa2750ec06df0b78a76ae5cafa0e0fd",
        "c7d976edeb9cd5406eae7aba4c05d6d92e81ae95"
    )

    $ws.Hyperlinks.Delete()

    for ($i = 0; $i -lt 4; $i++) {
        $row = $i + 2
        $xlfName = "$($names[$i].Substring(0, $names[$i].Length - 3)).$($xlfHashes[$i]).$xlfSuffix.xlf"

        $ws.Range("B$row").Value = ".md"
        $ws.Range("C$row").Value = "Ready for handoff"
        $ws.Range("E$row").Value = $handoffDatetime
        $ws.Range("H$row").Value = "0001-01-01 00:00:00"
        $ws.Range("J$row").Value = "Include"

        $ws.Hyperlinks.Add($ws.Range("A$row"), "$handoffLink$($names[$i])", $null, $null, $names[$i]) | Out-Null
        $ws.Hyperlinks.Add($ws.Range("D$row"), "$xlfLinkBase$xlfName", $null, $null, $xlfName) | Out-Null
    }

    # Reference Tokens ("Include" rows list the files that include them,
    # back-reference style) and Dependency From (for the caller rows):
    $ws.Range("K2").Value = "e2e\callerMd2.md`ne2e\callerMd1.md"
    $ws.Range("K3").Value = "e2e\callerMd1.md"
    $ws.Range("I4").Value = "e2e\calleeMd1.md`ne2e\calleeMd2.md"
    $ws.Range("I5").Value = "e2e\calleeMd1.md"
}

$e2eLink = "https://github.com/OpenLocalizationTest/oltest/blob/6cb7b01ae03b28ec25e20f4b9bba2f4dadeffb1d/e2e/"

$zhcn = $wb.Worksheets.Item("zh-cn")
Update-DetailSheet $zhcn "zh-cn" $e2eLink "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/4590bd4977577a8cd611f93e25dd72acc875c886/ol-handoff/OpenLocalizationTestOrg/oltest-zhcn-fly/yuwzho/ht/" "2016-03-19 17:14:59"

$dede = $wb.Worksheets.Item("de-de")
Update-DetailSheet $dede "de-de" $e2eLink "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/81ff8335ed19999995f076e2f8fdfe1263aa738d/ol-handoff/OpenLocalizationTestOrg/oltest-dede-fly/yuwzho/ht/" "2016-03-19 17:14:59"

# de-de sheet uses the overall "Latest Handoff Date" (17:15:08) for column E
# instead of the per-file handoff datetime (17:14:59) used on zh-cn.
for ($row = 2; $row -le 5; $row++) {
    $dede.Range("E$row").Value = "2016-03-19 17:15:08"
}

Write-Host "Report regenerated for handoff."
